# Daily attendance processing - 2025-11-02 20:21:52
#
# The "Recorded By" column (G) contains comma-separated lists of users.
# Normalize the position of the literal token "System" (capital S) within
# each list: it should immediately follow a leading lowercase "system"
# token when present, otherwise it should be moved to the very front of
# the list. All other tokens keep their relative order.
#
# Note: this runtime's PowerShell string comparison operators (-eq, -ne,
# -ceq, -cmatch, -like, ...) behave case-INsensitively, so a dedicated
# ordinal/char-code based comparison helper is used wherever the upper
# vs. lower case of "System" / "system" actually matters. Also note this
# runtime does not give functions their own variable scope, so the helper
# below uses variable names ($fnA, $fnB, $fnIdx, ...) that are never reused
# anywhere else in the script to avoid accidentally clobbering loop
# counters in the caller.

function Test-ExactEquals($fnA, $fnB) {
    if ($fnA.Length -ne $fnB.Length) { return $false }
    $fnCharsA = $fnA.ToCharArray()
    $fnCharsB = $fnB.ToCharArray()
    for ($fnIdx = 0; $fnIdx -lt $fnCharsA.Count; $fnIdx++) {
        $fnCodeA = [int]$fnCharsA[$fnIdx]
        $fnCodeB = [int]$fnCharsB[$fnIdx]
        if ($fnCodeA -ne $fnCodeB) { return $false }
    }
    return $true
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($rowNum = 2; $rowNum -le $lastRow; $rowNum++) {
    $cell = $ws.Cells.Item($rowNum, 7)
    $cellVal = $cell.Value2

    if ($cellVal -eq $null) { continue }
    if ($cellVal.GetType().Name -ne "String") { continue }
    if ($cellVal.IndexOf(",") -lt 0) { continue }

    $rawParts = $cellVal.Split(",")

    $trimmedParts = @()
    foreach ($rawPart in $rawParts) {
        $trimmedParts += $rawPart.Trim()
    }

    # locate an exact (case-sensitive) "System" token
    $systemPos = -1
    for ($scanPos = 0; $scanPos -lt $trimmedParts.Count; $scanPos++) {
        if (Test-ExactEquals $trimmedParts[$scanPos] "System") {
            $systemPos = $scanPos
        }
    }

    if ($systemPos -lt 0) { continue }

    # build the list without the "System" token
    $remainingParts = @()
    for ($copyPos = 0; $copyPos -lt $trimmedParts.Count; $copyPos++) {
        if ($copyPos -ne $systemPos) { $remainingParts += $trimmedParts[$copyPos] }
    }

    $leadsWithLowerSystem = ($remainingParts.Count -gt 0) -and (Test-ExactEquals $remainingParts[0] "system")

    $finalParts = @()
    if ($leadsWithLowerSystem) {
        $finalParts += $remainingParts[0]
        $finalParts += "System"
        for ($appendPos = 1; $appendPos -lt $remainingParts.Count; $appendPos++) { $finalParts += $remainingParts[$appendPos] }
    } else {
        $finalParts += "System"
        for ($appendPos = 0; $appendPos -lt $remainingParts.Count; $appendPos++) { $finalParts += $remainingParts[$appendPos] }
    }

    $newCellVal = [string]::Join(", ", $finalParts)

    if (-not (Test-ExactEquals $newCellVal $cellVal)) {
        $cell.Value = $newCellVal
    }
}
